# Generate Report for Handoff
# Updates the localization-status report: the "In Translation" rows have
# now been handed off, so the Status columns move to "Ready for handoff"
# and the handoff timestamps are refreshed to the moment the report was
# (re)generated.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---------------------------------------------------
# B2 = zh-cn status, C2 = de-de status, D2 = latest handoff date/time
$overview.Range("B2").Value = "Ready for handoff"
$overview.Range("C2").Value = "Ready for handoff"
$overview.Range("D2").Value = "2016-29-20 08:29:01"

# --- zh-cn sheet --------------------------------------------------------
# C2 = Status, E2 = Latest Handoff Datetime
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("E2").Value = "2016-03-20 08:28:58"

# --- de-de sheet --------------------------------------------------------
# C2 = Status, E2 = Latest Handoff Datetime
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("E2").Value = "2016-03-20 08:29:01"
